# Auto-generated edit script: updates crafting-profit cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 242.61111
$ws.Range("I9").Value = 103.76923
$ws.Range("J9").Value = 603.6
$ws.Range("K9").Value = 103.76923
$ws.Range("L9").Value = 603.6
$ws.Range("M9").Value = 65.23077000000001
$ws.Range("N9").Value = -941.6

$ws.Range("H48").Value = 4000
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12584

$ws.Range("H56").Value = 4000
$ws.Range("J56").Value = 4000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13068

$ws.Range("H99").Value = 478.14285
$ws.Range("I99").Value = 282.66666
$ws.Range("J99").Value = 624.75
$ws.Range("K99").Value = 847.9999799999999
$ws.Range("L99").Value = 1874.25
$ws.Range("M99").Value = 650.0000200000001
$ws.Range("N99").Value = -4870.25

$ws.Range("H132").Value = 3025.3547
$ws.Range("I132").Value = 2056.04
$ws.Range("K132").Value = 6168.12
$ws.Range("M132").Value = -3638.12

$ws.Range("H138").Value = 2886.8333
$ws.Range("J138").Value = 3288
$ws.Range("L138").Value = 9864
$ws.Range("N138").Value = -20144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2684.8386
$ws.Range("I32").Value = 2684.8386
$ws.Range("K32").Value = 2684.8386
$ws.Range("M32").Value = -2397.8386

$ws.Range("H80").Value = 152527.5
$ws.Range("J80").Value = 176703.33
$ws.Range("L80").Value = 176703.33
$ws.Range("N80").Value = -178699.33

$ws.Range("H83").Value = 152527.5
$ws.Range("J83").Value = 176703.33
$ws.Range("L83").Value = 530109.99
$ws.Range("N83").Value = -540093.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2628.5
$ws.Range("I20").Value = 2628.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2628.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2381.5
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 314.92307
$ws.Range("I22").Value = 317.63635
$ws.Range("K22").Value = 317.63635
$ws.Range("M22").Value = -144.63635

$ws.Range("H33").Value = 3190.3333
$ws.Range("I33").Value = 1828.4
$ws.Range("K33").Value = 1828.4
$ws.Range("M33").Value = -1492.4

$ws.Range("H99").Value = 1692
$ws.Range("I99").Value = 1462.3636
$ws.Range("K99").Value = 1462.3636
$ws.Range("M99").Value = 35.63640000000009

$ws.Range("H134").Value = 3998.5334
$ws.Range("I134").Value = 4036.7693
$ws.Range("K134").Value = 12110.3079
$ws.Range("M134").Value = -9575.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2129.75
$ws.Range("I31").Value = 1834.6364
$ws.Range("K31").Value = 1834.6364
$ws.Range("M31").Value = -1539.6364

$ws.Range("H34").Value = 2129.75
$ws.Range("I34").Value = 1834.6364
$ws.Range("K34").Value = 1834.6364
$ws.Range("M34").Value = -1632.6364

$ws.Range("H134").Value = 912.4643
$ws.Range("I134").Value = 912.4643
$ws.Range("K134").Value = 2737.3929
$ws.Range("M134").Value = -202.3928999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1520.84
$ws.Range("I12").Value = 944.1111
$ws.Range("J12").Value = 1845.25
$ws.Range("K12").Value = 2832.3333
$ws.Range("L12").Value = 5535.75
$ws.Range("M12").Value = -2659.3333
$ws.Range("N12").Value = -5881.75

$ws.Range("H37").Value = 69983
$ws.Range("J37").Value = 69983
$ws.Range("L37").Value = 209949
$ws.Range("N37").Value = -210173

$ws.Range("H47").Value = 1118.091
$ws.Range("I47").Value = 1118.091
$ws.Range("K47").Value = 3354.273
$ws.Range("M47").Value = -2923.273

$ws.Range("H131").Value = 590590.1
$ws.Range("I131").Value = 1033
$ws.Range("J131").Value = 716923.8
$ws.Range("K131").Value = 3099
$ws.Range("L131").Value = 2150771.4
$ws.Range("M131").Value = 1941
$ws.Range("N131").Value = -2160851.4

$ws.Range("H141").Value = 12618.2
$ws.Range("I141").Value = 13139.5
$ws.Range("K141").Value = 39418.5
$ws.Range("M141").Value = -34238.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 362029.75
$ws.Range("I7").Value = 480039.66
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 480039.66
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -479927.66
$ws.Range("N7").Value = -8224

$ws.Range("H8").Value = 362029.75
$ws.Range("I8").Value = 480039.66
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 480039.66
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -479900.66
$ws.Range("N8").Value = -8278

$ws.Range("H24").Value = 10000
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("N24").Value = -10346

$ws.Range("H70").Value = 9977
$ws.Range("J70").Value = 9977
$ws.Range("L70").Value = 9977
$ws.Range("N70").Value = -10517

$ws.Range("H73").Value = 9977
$ws.Range("J73").Value = 9977
$ws.Range("L73").Value = 9977
$ws.Range("N73").Value = -11849

$ws.Range("H95").Value = 26899
$ws.Range("J95").Value = 26899
$ws.Range("L95").Value = 26899
$ws.Range("N95").Value = -32391

$ws.Range("H113").Value = 1606.5
$ws.Range("I113").Value = 1327.8
$ws.Range("K113").Value = 1327.8
$ws.Range("M113").Value = 842.2

$ws.Range("H126").Value = 6814.3335
$ws.Range("I126").Value = 7974.75
$ws.Range("J126").Value = 4493.5
$ws.Range("K126").Value = 23924.25
$ws.Range("L126").Value = 13480.5
$ws.Range("M126").Value = -21454.25
$ws.Range("N126").Value = -18420.5

$ws.Range("H132").Value = 1216.2727
$ws.Range("I132").Value = 1216.2727
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3648.8181
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1118.8181
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4952.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 4952.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 4952.5
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -5404.5

$ws.Range("H55").Value = 240.42857
$ws.Range("J55").Value = 281
$ws.Range("L55").Value = 281
$ws.Range("N55").Value = -627

$ws.Range("H68").Value = 2840.7273
$ws.Range("I68").Value = 2805.5557
$ws.Range("J68").Value = 2999
$ws.Range("K68").Value = 2805.5557
$ws.Range("L68").Value = 2999
$ws.Range("M68").Value = -2056.5557
$ws.Range("N68").Value = -4497

$ws.Range("H71").Value = 2840.7273
$ws.Range("I71").Value = 2805.5557
$ws.Range("J71").Value = 2999
$ws.Range("K71").Value = 14027.7785
$ws.Range("L71").Value = 14995
$ws.Range("M71").Value = -10283.7785
$ws.Range("N71").Value = -22483

$ws.Range("H101").Value = 15907.667
$ws.Range("J101").Value = 15907.667
$ws.Range("L101").Value = 15907.667
$ws.Range("N101").Value = -22397.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I21").Value = 3000
$ws.Range("J21").Value = 7500
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 7500
$ws.Range("M21").Value = -2765
$ws.Range("N21").Value = -7970

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H30").Value = 9500
$ws.Range("I30").Value = 15000
$ws.Range("J30").Value = 4000
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 4000
$ws.Range("M30").Value = -14893
$ws.Range("N30").Value = -4214

$ws.Range("I35").Value = 3000
$ws.Range("J35").Value = 7500
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 7500
$ws.Range("M35").Value = -2710
$ws.Range("N35").Value = -8080

$ws.Range("H62").Value = 34000.5
$ws.Range("I62").Value = 40002
$ws.Range("J62").Value = 32000
$ws.Range("K62").Value = 40002
$ws.Range("L62").Value = 32000
$ws.Range("M62").Value = -39378
$ws.Range("N62").Value = -33248

$ws.Range("H65").Value = 34000.5
$ws.Range("I65").Value = 40002
$ws.Range("J65").Value = 32000
$ws.Range("K65").Value = 200010
$ws.Range("L65").Value = 160000
$ws.Range("M65").Value = -196890
$ws.Range("N65").Value = -166240

$ws.Range("H103").Value = 20000
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws.Range("H132").Value = 5086.0454
$ws.Range("I132").Value = 4550.0557
$ws.Range("K132").Value = 13650.1671
$ws.Range("M132").Value = -11120.1671

